# Update industrial biomass extension.
# Adds a new "other_biomass" worksheet (matching the layout used by the
# existing percent-extension sheets) after the last existing sheet, fills
# in its sector/X1850 data, and fixes up the selection on the "petroleum"
# sheet that the author left behind while working on this sheet.

$wb = $excel.ActiveWorkbook

# --- petroleum sheet: selection left on the header row when the author
#     moved on to the new sheet ---
$petroleum = $wb.Worksheets.Item("petroleum")
$petroleum.Range("A1:C1").Select()

# --- Add the new sheet at the end of the workbook, making it active ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "other_biomass"

# Column layout: A is narrower than the default to fit the long sector
# names, like the other extension sheets in this workbook.
$ws.Columns.Item(1).ColumnWidth = 24.6640625

$data = @(
    @("sector", "X1850"),
    @("1A1a_Electricity-autoproducer", 0),
    @("1A1a_Electricity-public", 0),
    @("1A1a_Heat-production", 0),
    @("1A3ai_International-aviation", 0),
    @("1A3aii_Domestic-aviation", 0),
    @("1A3b_Road", 0),
    @("1A3c_Rail", 0),
    @("1A3dii_Domestic-navigation", 0),
    @("1A3eii_Other-transp", 0),
    @("1A4a_Commercial-institutional", 0),
    @("1A4c_Agriculture-forestry-fishing", 0),
    @("1A5_Other-unspecified", 1)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

$ws.Range("B1").Select()
